# Automatische test-sync: 2025-06-22 21:41:50
# Adds new mail-log entry (row 44) to the "Logs" sheet and refreshes the
# "Dashboard" summary counts to reflect it.

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append the new e-mail entry -----------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A44").Value = "Korting voor wederverkopers?"
$logs.Range("B44").Value = "mailmind.test@zohomail.eu"
$logs.Range("C44").Value = "Biedt u speciale prijzen voor wederverkopers?"
$logs.Range("D44").Value = "Offerte / Prijsaanvraag"
$logs.Range("E44").Value = "Beste klant,`nDank u voor uw interesse in onze producten. Ja, wij bieden speciale prijzen aan voor wederverkopers. Om hier meer informatie over te krijgen en om te weten te komen hoe wij u verder kunnen helpen, kunt u het beste contact opnemen met ons verkoopteam via verkoop@bedrijfsnaam.nl.`nMet vriendelijke groet,`n[Bedrijfsnaam] Team"
$logs.Range("F44").Value = "2025-06-22 21:41:26"
$logs.Range("G44").Value = "Ja"

# Re-fit the row height so the multi-line "Antwoord" text doesn't leave a
# stray custom row height behind (matches the other rows in the sheet).
$logs.Rows.Item(44).AutoFit()

# Conditional formatting ranges on columns D and G need to grow to include
# the new row 44.
$dFormats = $logs.Range("D2").FormatConditions
for ($i = 1; $i -le $dFormats.Count; $i++) {
    $dFormats.Item($i).ModifyAppliesToRange($logs.Range("D2:D44"))
}

$gFormats = $logs.Range("G2").FormatConditions
for ($i = 1; $i -le $gFormats.Count; $i++) {
    $gFormats.Item($i).ModifyAppliesToRange($logs.Range("G2:G44"))
}

# --- Dashboard sheet: update category counts ------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

# "Offerte / Prijsaanvraag" moves up to rank 4 (count 5), pushing
# "Productinformatie" down to rank 8 (count 4).
$dash.Range("A4").Value = "Offerte / Prijsaanvraag"
$dash.Range("A8").Value = "Productinformatie"
$dash.Range("B8").Value = 4
